$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab title from "Through 2022-12-18" to "Through 2022-12-19"
$ws.Name = "Through 2022-12-19"

# Update label for the December row
$ws.Range("A13").Value = "December (through 12-19)"

# Update December monthly figures (row 13)
$ws.Range("B13").Value = 24
$ws.Range("C13").Value = 60
$ws.Range("D13").Value = 74
$ws.Range("E13").Value = 42
$ws.Range("F13").Value = 30
$ws.Range("G13").Value = 89
$ws.Range("H13").Value = 135
$ws.Range("I13").Value = 81

# Update Total figures (row 14)
$ws.Range("B14").Value = 315
$ws.Range("C14").Value = 623
$ws.Range("D14").Value = 895
$ws.Range("E14").Value = 724
$ws.Range("F14").Value = 564
$ws.Range("G14").Value = 1353
$ws.Range("H14").Value = 1778
$ws.Range("I14").Value = 1598
